$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @('B','C','D','E','F','I','J','K','L','M','N')

$data = @(
    @(2,1.02,1.011267060624746,1.040897159324743,1.013467924374028,1.042723146085332,1.035189460091314,1.016516520928296,1.043678206412957,1.01632898695471,1.045499028740368,1.009642935836688),
    @(3,1.02,1.012187912163178,1.041515364982076,1.014247300626775,1.043570948546279,1.035308032181936,1.017070705920199,1.044107434926798,1.016913169610881,1.046157624116618,1.009828002616436),
    @(4,1.02,1.012784202276114,1.041913132728877,1.014752391356762,1.044117744335275,1.035381933978818,1.017429160808853,1.044382235122282,1.017291303924076,1.0465813519545,1.009947644391614),
    @(5,1.02,1.01303498636883,1.042079812506646,1.014964917326041,1.044347186501278,1.035412325264072,1.017579820734225,1.044497055195688,1.017450301392447,1.046758903578535,1.009997915361388),
    @(6,1.02,1.013077100174195,1.042107766930905,1.015000612226713,1.044385685516441,1.035417388352756,1.017605115144946,1.044516292529163,1.017476999482659,1.046788680985669,1.010006354517187),
    @(7,1.02,1.012787552859157,1.04191536204433,1.01475523041081,1.04412081184866,1.035382342730727,1.017431174071062,1.044383772130488,1.017293428340322,1.046583726704177,1.009948316219123),
    @(8,1.02,1.011578175488351,1.041106549761852,1.013731155415298,1.04301003386746,1.035230115638594,1.016703838683976,1.04382387347728,1.016526386835833,1.045722105293099,1.009705502320287),
    @(9,1.02,1.009450494622893,1.03966417436083,1.011932657729253,1.041039105631031,1.034940315948077,1.015421154602985,1.042814845279282,1.015175798308388,1.044185332349406,1.009276816700478),
    @(10,1.02,1.008034383486343,1.038691205394766,1.010737812775054,1.03971615272959,1.034732700684175,1.014565401388698,1.042127216085648,1.014276166734804,1.043148526260373,1.008990499629935),
    @(11,1.02,1.007421759311603,1.038267227520106,1.010221434214447,1.039141193051929,1.034639397565506,1.014194713444786,1.041825946385221,1.013886809935324,1.042696690610843,1.008866400855336),
    @(12,1.02,1.007194288642084,1.038109343464763,1.010029779416815,1.03892731224644,1.034604230361081,1.014057003174901,1.041713514278104,1.013742215000496,1.042528426426434,1.008820287134314),
    @(13,1.02,1.00724307800181,1.038143228186845,1.010070883157752,1.038973204633383,1.034611796917826,1.014086543397093,1.041737655169085,1.013773229752692,1.042564539193859,1.008830179488789),
    @(14,1.02,1.00740295477133,1.03825418491782,1.0102055888769,1.039123520015008,1.034636501036016,1.014183330678834,1.041816663454142,1.01387485706484,1.04268279064751,1.008862589444084),
    @(15,1.02,1.007501471445737,1.038322496110932,1.010288605591673,1.039216092582402,1.034651654477777,1.014242961860632,1.041865273275066,1.01393747694792,1.042755592049688,1.008882555924523),
    @(16,1.02,1.008075053237709,1.038719287223152,1.01077210424664,1.03975426659177,1.034738821281694,1.014589999833488,1.042147136334182,1.014302011146171,1.043178452362526,1.00899873312562),
    @(17,1.02,1.008434996850964,1.038967468979315,1.011075658127952,1.040091284944366,1.03479258812651,1.014807650462859,1.042323000125855,1.01453072531067,1.04344292878708,1.009071575713523),
    @(18,1.02,1.008644999967296,1.039111970867094,1.011252812093672,1.040287658093716,1.034823620869615,1.014934588710237,1.042425238575125,1.01466414863902,1.043596914255089,1.009114051806782),
    @(19,1.02,1.008716614751333,1.039161198390617,1.011313233280017,1.040354581602051,1.034834146468999,1.014977869022622,1.042460041526432,1.014709645605992,1.043649371904657,1.009128533051995),
    @(20,1.02,1.008396372705717,1.0389408681298,1.011043079716915,1.040055147144596,1.034786853427409,1.014784300018508,1.042304166750947,1.014506184552686,1.043414581828345,1.009063761606582),
    @(21,1.02,1.007355872687672,1.038221521944703,1.010165917215427,1.039079264564919,1.034629240365505,1.014154829788517,1.041793412019058,1.013844929531361,1.042647980458461,1.008853046011609),
    @(22,1.02,1.006702162563616,1.037766928427044,1.009615285759563,1.038463865959604,1.034527191072961,1.013758940100667,1.041469231402491,1.013429343533801,1.04216348759418,1.008720457303625),
    @(23,1.02,1.007048659588376,1.038008135472702,1.009907102475263,1.038790272440601,1.034581568730473,1.01396881946963,1.041641374028942,1.013649636956625,1.042420562707833,1.008790754773527),
    @(24,1.02,1.008413825127269,1.038952888706978,1.011057800201662,1.040071476891016,1.034789445706286,1.014794851121534,1.042312677792777,1.014517273408636,1.043427391461423,1.009067292501504),
    @(25,1.02,1.010000142510728,1.040039079932728,1.012396886501438,1.04155023435521,1.035017782516342,1.015752874638931,1.04307834631716,1.015524828371555,1.044584801269465,1.009387736660142)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $entry[$i + 1]
        $ws.Range("$col$row").Value = $val
    }
}
